$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17 (A17 = 5) ---
$ws.Range("M17").Value2 = $false
$ws.Range("N17").Value2 = "Missing 1 solid, Solid counted, mass found by matching to expected solution results (non-detect/detect)"

# --- Row 18 (A18 = 10) ---
$ws.Range("C18").Value2 = -0.024383656501640903
$ws.Range("D18").Value2 = 0.0035532598955532839
$ws.Range("E18").Value2 = 261.2389380530974
$ws.Range("F18").Value2 = 11.844199651844281
$ws.Range("H18").Value2 = 1.4516493526175951
$ws.Range("I18").Value2 = 0.065815704523261576
$ws.Range("L18").Value2 = 42517
$ws.Range("M18").Value2 = $false
$ws.Range("N18").Value2 = "Total activity exceeds possible initial activity. Solid counted, mass found by matching to expected solution results (non-detect/detect)."

# --- Row 19 (A19 = 50) ---
$ws.Range("C19").Value2 = -0.010957117121438813
$ws.Range("D19").Value2 = 0.0042526244585913853
$ws.Range("E19").Value2 = 937.46312684365785
$ws.Range("F19").Value2 = 14.175414861971241
$ws.Range("H19").Value2 = 1.0405396007902779
$ws.Range("I19").Value2 = 0.015734038064168072
$ws.Range("L19").Value2 = 42517

# --- Row 20 (A20 = 100) ---
$ws.Range("C20").Value2 = 0.017402373036688297
$ws.Range("D20").Value2 = 0.012302237781342635
$ws.Range("E20").Value2 = 1741.5929203539824
$ws.Range("F20").Value2 = 41.007459271142103
$ws.Range("H20").Value2 = 0.96776623507839676
$ws.Range("I20").Value2 = 0.022786975076183512
$ws.Range("L20").Value2 = 42517

# --- Row 21 (A21 = 500) ---
$ws.Range("E21").Value2 = 8520.0337221815589
$ws.Range("F21").Value2 = 357.94377699140438
$ws.Range("H21").Value2 = 0.94448928246606678
$ws.Range("I21").Value2 = 0.039679897065858227
$ws.Range("L21").Value2 = 42517

# --- Selection change on the active pane ---
$ws.Range("I14").Select()
